$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New records (rows 7-21) appended below the existing 5 records.
# Columns: A=Tipo Documento, B=Numero Documento, C=Nombres y Apellidos,
#          D=Dia, E=Mes, F=Año
$newRecords = @(
    ,@("CC", 38253096, "NA", 5, "NOVIEMBRE", 1980)
    ,@("CC", 65769149, "NA", 20, "JUNIO", 1994)
    ,@("CC", 1110574193, "NA", 2, "SEPTIEMBRE", 2014)
    ,@("CC", 1110549747, "NA", 9, "AGOSTO", 2010)
    ,@("CC", 1617567, "NA", 28, "ABRIL", 1956)
    ,@("CC", 65735017, "NA", 2, "SEPTIEMBRE", 1985)
    ,@("CC", 12185483, "NA", 6, "JULIO", 1972)
    ,@("CC", 41922408, "NA", 31, "JULIO", 1989)
    ,@("CC", 28612802, "NA", 19, "ABRIL", 1989)
    ,@("CC", 28603778, "NA", 21, "ENERO", 1974)
    ,@("CC", 65733352, "NA", 22, "FEBRERO", 1985)
    ,@("CC", 28604216, "NA", 25, "JULIO", 1975)
    ,@("CC", 38257409, "NA", 5, "NOVIEMBRE", 1981)
    ,@("CC", 65738644, "NA", 5, "NOVIEMBRE", 1986)
    ,@("CC", 38261159, "NA", 18, "NOVIEMBRE", 1982)
)

$startRow = 7
for ($i = 0; $i -lt $newRecords.Count; $i++) {
    $r = $startRow + $i
    $rec = $newRecords[$i]
    $ws.Cells.Item($r, 1).Value = $rec[0]
    $ws.Cells.Item($r, 2).Value = $rec[1]
    $ws.Cells.Item($r, 3).Value = $rec[2]
    $ws.Cells.Item($r, 4).Value = $rec[3]
    $ws.Cells.Item($r, 5).Value = $rec[4]
    $ws.Cells.Item($r, 6).Value = $rec[5]
}

# Update the saved selection / active cell to reflect where the user
# was working after entering the new data (matches the workbook state).
$ws.Range("H15").Select()
